$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old row 10 (B10 had placeholder "x")
$ws.Range("B10").Value = $null

# Add new data rows 8 and 9
$ws.Range("A8").Value = "2018-01-02 22:08:12.41Z"
$ws.Range("B8").Value = "input date is not in correct format"

$ws.Range("A9").Value = "2018-01-0"
$ws.Range("B9").Value = "input date is not in correct format"

# Match style of A4/A5/A7 (left aligned) for the new A8/A9 cells
$ws.Range("A8").HorizontalAlignment = -4131
$ws.Range("A9").HorizontalAlignment = -4131

# Update selection to match target (activeCell B9, sqref B9)
$ws.Range("B9").Select()
